$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "runs"/"balls" columns (C/D) store numeric-looking data as TEXT
# (see the ignoredErrors/numberStoredAsText marker on this sheet). Force
# each destination cell to Text format first so Excel keeps writing the
# new value as text instead of silently converting it to a number.
$ws.Range("C2:D2").NumberFormat = "@"
$ws.Range("C3:D3").NumberFormat = "@"
$ws.Range("C4:D4").NumberFormat = "@"
$ws.Range("C6:D6").NumberFormat = "@"

$ws.Range("C2").Value = "3"
$ws.Range("D2").Value = "3"

$ws.Range("C3").Value = "6"
$ws.Range("D3").Value = "13"

$ws.Range("C4").Value = "0"
$ws.Range("D4").Value = "2"

$ws.Range("C6").Value = "1"
$ws.Range("D6").Value = "1"
